$d = $word.ActiveDocument

# Helper wildcard-off Find & Replace all occurrences
function Replace-Text($findText, $replaceText) {
    $d.Content.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, $replaceText, 2) | Out-Null
}

# 1) ". For futures challenges" -> ". For future challenges"
Replace-Text "futures challenges" "future challenges"

# 2) Paragraph "I'm looking forward..." restructuring:
#    - "evaluate the all given" -> "evaluate all given"
#    - "...this seems to me more like a random choice than a math based evaluate. Normally you are given a 80% of a given dataset, train a model, and compare to the rest 20%."
#      -> "...this seems to me the evaluation would be difficult. "
Replace-Text "evaluate the all given" "evaluate all given"
Replace-Text "this seems to me more like a random choice than a math based evaluate. Normally you are given a 80% of a given dataset, train a model, and compare to the rest 20%." "this seems to me the evaluation would be difficult. "

# 3) "I did not uses tripadvisor" text unchanged; nothing to do for wording (run-merge cosmetic).

# 4) "Because of the small dataset " -> "Because of the small size of the dataset "
Replace-Text "Because of the small dataset " "Because of the small size of the dataset "

# 5) "requested most often and how what kind of service" -> "requested most often and what kind of service"
Replace-Text "requested most often and how what kind of service" "requested most often and what kind of service"

# 6) After "I just gave it my best guess." add two empty paragraphs, then bookmark stays,
#    then "future: " run with lastRenderedPageBreak, then a new paragraph with the final text.

$p18 = $d.Paragraphs(18)
$insertAt = $p18.Range.End - 1
$ip = $d.Range($insertAt, $insertAt)
$ip.Text = "`r`r"

# Now paragraph 19 should be a new empty paragraph, paragraph 20 the other new empty
# paragraph that originally held the bookmark. Re-fetch paragraph objects fresh.

$p20 = $d.Paragraphs(20)
$p20Start = $p20.Range.Start
$insertText = $d.Range($p20Start, $p20Start)
$insertText.InsertAfter("future: `r")

$p21 = $d.Paragraphs(21)
$p21End = $p21.Range.End - 1
$insertText2 = $d.Range($p21End, $p21End)
$insertText2.InsertAfter("The designed model is based on a cold start concept. There are not ratings or any user input involved. A future design should use user input data (for example just thump up / down) do reevaluate the LSA Model.")

Write-Host "Edits applied"
